$d = $word.ActiveDocument

$old = "Definisanje scenarija upotrebe pri registrovanju administratora. "
$new = "Definisanje scenarija upotrebe pri registrovanju administratora, sa grafičkim opisom priloženim u prototipu koji se nalazi u drugom folderu. "

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Replace result: $found"
